# feat: treemap para recursos de biomasa
#
# The "Recursos" sheet is duplicated into a new "Recursos_old" sheet
# (placed right after "Recursos", before "Datos"), preserving the old
# figures. The live "Recursos" sheet is then updated with new data.

$wb = $excel.ActiveWorkbook

$wsRecursos = $wb.Worksheets.Item("Recursos")

# Duplicate "Recursos" -> new copy placed immediately after it; rename to
# "Recursos_old" so it sits between "Recursos" and "Datos". Grab the new
# sheet by position (right after "Recursos") rather than by its default
# "Recursos (2)" name, since that auto-generated name isn't guaranteed.
$origIndex = $wsRecursos.Index
$wsRecursos.Copy($null, $wsRecursos)
$wsOld = $wb.Worksheets.Item($origIndex + 1)
$wsOld.Name = "Recursos_old"

# Update the values on the live "Recursos" sheet.
$wsRecursos.Range("B2").Value = 2563
$wsRecursos.Range("B3").Value = 4600
$wsRecursos.Range("B4").Value = 151000
$wsRecursos.Range("B7").Value = 0
$wsRecursos.Range("B8").Value = 28128
$wsRecursos.Range("F8").Value = 0.25

# Restore the active sheet / selection so "Recursos" stays the active tab
# with B8 selected (and "Recursos_old" keeps the original G11 selection).
$wsRecursos.Activate() | Out-Null
$wsRecursos.Range("B8").Select() | Out-Null
